$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 44; Date = "2025-06-25 17:50:26"; Score = 1; Episode = 1 },
    @{ Row = 45; Date = "2025-06-25 17:50:31"; Score = 0; Episode = 2 },
    @{ Row = 46; Date = "2025-06-25 17:50:35"; Score = 1; Episode = 3 },
    @{ Row = 47; Date = "2025-06-25 17:50:40"; Score = 0; Episode = 4 },
    @{ Row = 48; Date = "2025-06-25 17:50:44"; Score = 0; Episode = 5 },
    @{ Row = 49; Date = "2025-06-25 17:50:54"; Score = 0; Episode = 6 },
    @{ Row = 50; Date = "2025-06-25 17:51:01"; Score = 0; Episode = 7 },
    @{ Row = 51; Date = "2025-06-25 17:51:07"; Score = 1; Episode = 8 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "Q Learning"
    $ws.Cells.Item($row, 3).Value = "MontyHall LV2"
    $ws.Cells.Item($row, 4).Value = $r.Score
    $ws.Cells.Item($row, 5).Value = $r.Episode
}
